# Updated cryptos list - applies latest price/volume(1h) refresh to the
# "cryptos" worksheet, matching the upstream GitHub Actions scrape.
#
# Column D ("Price") cells are written as text (values that look like plain
# numbers are prefixed with a leading apostrophe) so Excel keeps them as
# strings instead of silently converting them to numeric values, exactly
# like the source data. Column E ("Volume(1h)") values already contain a
# "%" sign and padding spaces, so they naturally stay text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range('D2').Value = '45.444.18'
$ws.Range('E2').Value = '  +7.54%  '

# Row 3
$ws.Range('D3').Value = '2.381.80'
$ws.Range('E3').Value = '  +4.61%  '

# Row 4
$ws.Range('E4').Value = '  -1.44%  '

# Row 5
$ws.Range('D5').Value = '''113.24'
$ws.Range('E5').Value = '  +10.70%  '

# Row 6
$ws.Range('D6').Value = '''318.46'
$ws.Range('E6').Value = '  +3.07%  '

# Row 7
$ws.Range('D7').Value = '''0.636'

# Row 8
$ws.Range('E8').Value = '  -0.29%  '

# Row 9
$ws.Range('D9').Value = '''0.628'
$ws.Range('E9').Value = '  +5.36%  '

# Row 10
$ws.Range('D10').Value = '''42.74'
$ws.Range('E10').Value = '  +11.49%  '

# Row 11
$ws.Range('D11').Value = '''0.0932'
$ws.Range('E11').Value = '  +4.10%  '

# Row 12
$ws.Range('D12').Value = '''8.69'
$ws.Range('E12').Value = '  +6.49%  '

# Row 13
$ws.Range('E13').Value = '  +5.58%  '

# Row 14
$ws.Range('E14').Value = '  +1.77%  '

# Row 15
$ws.Range('D15').Value = '''15.86'
$ws.Range('E15').Value = '  +5.44%  '

# Row 16
$ws.Range('D16').Value = '2.747.04'
$ws.Range('E16').Value = '  +4.69%  '

# Row 17
$ws.Range('D17').Value = '2.380.98'
$ws.Range('E17').Value = '  +4.47%  '

# Row 18
$ws.Range('D18').Value = '45.343.46'
$ws.Range('E18').Value = '  +6.50%  '

# Row 19
$ws.Range('D19').Value = '''7.61'
$ws.Range('E19').Value = '  +5.59%  '

# Row 20
$ws.Range('E20').Value = '  +4.20%  '

# Row 21
$ws.Range('D21').Value = '''13.31'
$ws.Range('E21').Value = '  +2.11%  '

# Row 22
$ws.Range('D22').Value = '''74.89'
$ws.Range('E22').Value = '  +3.47%  '

# Row 23
$ws.Range('D23').Value = '''3.55'
$ws.Range('E23').Value = '  +5.65%  '

# Row 24
$ws.Range('D24').Value = '''270.37'
$ws.Range('E24').Value = '  +3.88%  '

# Row 25
$ws.Range('D25').Value = '''2.38'
$ws.Range('E25').Value = '  +10.45%  '

# Row 26
$ws.Range('E26').Value = '  -0.42%  '

# Row 27
$ws.Range('D27').Value = '''11.27'
$ws.Range('E27').Value = '  +6.37%  '

# Row 28
$ws.Range('E28').Value = '  +10.07%  '

# Row 29
$ws.Range('E29').Value = '  -0.28%  '

# Row 30
$ws.Range('D30').Value = '''39.52'
$ws.Range('E30').Value = '  +11.11%  '

# Row 31
$ws.Range('D31').Value = '''23.01'
$ws.Range('E31').Value = '  +4.97%  '

# Row 32
$ws.Range('D32').Value = '''0.0942'
$ws.Range('E32').Value = '  +11.76%  '

# Row 33
$ws.Range('D33').Value = '''170.01'
$ws.Range('E33').Value = '  +3.46%  '

# Row 34
$ws.Range('D34').Value = '''2.97'
$ws.Range('E34').Value = '  +16.68%  '

# Row 35
$ws.Range('E35').Value = '  +3.11%  '

# Row 36
$ws.Range('D36').Value = '''4.94'
$ws.Range('E36').Value = '  +10.91%  '

# Row 37
$ws.Range('E37').Value = '  +8.21%  '

# Row 38
$ws.Range('D38').Value = '''3.06'
$ws.Range('E38').Value = '  +14.38%  '

# Row 39
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = '''4.01'
$ws.Range('E39').Value = '  +10.79%  '

# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.0364'
$ws.Range('E40').Value = '  +5.46%  '

# Row 41
$ws.Range('D41').Value = '''1.73'
$ws.Range('E41').Value = '  +12.37%  '

# Row 42
$ws.Range('D42').Value = '''106.31'
$ws.Range('E42').Value = '  +8.67%  '

# Row 43
$ws.Range('E43').Value = '  +7.93%  '

# Row 44
$ws.Range('D44').Value = '''13.47'
$ws.Range('E44').Value = '  +14.11%  '

# Row 45
$ws.Range('D45').Value = '''71.65'
$ws.Range('E45').Value = '  +5.10%  '

# Row 46
$ws.Range('E46').Value = '  -0.55%  '

# Row 47
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').Value = '''5.79'
$ws.Range('E47').Value = '  +14.02%  '

# Row 48
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '''116.75'
$ws.Range('E48').Value = '  +6.81%  '

# Row 49
$ws.Range('D49').Value = '''1.65'
$ws.Range('E49').Value = '  +21.55%  '

# Row 50
$ws.Range('D50').Value = '''9.32'
$ws.Range('E50').Value = '  +8.50%  '

# Row 51
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').Value = '''78.64'
$ws.Range('E51').Value = '  +3.40%  '
